# ANOR_suppl_mat.xlsx edit script
# - rename "validation_table" -> "stereology_visualization" (Excel auto-updates
#   the _xlchart defined names that reference the sheet by name)
# - update a few header / formula-description cells on the "Equations" sheet
# - nudge the saved selections / scroll position on a couple of sheets to
#   match the authored view state

$wb = $excel.ActiveWorkbook

# 1. Rename the sheet. Renaming via the Name property keeps every
#    definedName / formula reference that points at the old name in sync.
$validation = $wb.Worksheets.Item("validation_table")
$validation.Name = "stereology_visualization"

# 2. Content edits on the "Equations" sheet (order matches the authored
#    edit sequence so the shared-string table is appended in the same
#    order as the reference workbook).
$eq = $wb.Worksheets.Item("Equations")

$eq.Range("C14").Value = "(x17)/(x22+x15) = 1"
$eq.Range("C2").Value = "(x19+x20)/(2*(x10+x18)) = 1"
$eq.Range("B1").Value = "Percent-error equation"
$eq.Range("A1").Value = "Weights"
$eq.Range("C1").Value = "Measurement equation"

# 3. View-state tweaks matching the authored workbook.
$eq.Application.ActiveWindow.ScrollRow = 1
$eq.Range("C2").Select()

$validation.Select()
$validation.Range("R5").Select()

$eq.Select()
